$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.211459279060364
$ws.Range("B1").Value = 2.269507169723511
$ws.Range("C1").Value = 3.399822950363159
$ws.Range("D1").Value = 2.420431613922119
$ws.Range("E1").Value = 1.313086748123169
